# Populate the Help Items sheet with card-statistics content (rows 4-20)
# and the associated cell formatting (wrap text + row heights).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sayfa1")

# Row 4
$b4 = @'
BANKER
'@
$c4 = @'
A Banker who plays in the game must keep hislher personal funds separate from those of the Bank.
'@
$ws.Range("B4").Value = $b4
$ws.Range("C4").Value = $c4

# Row 5
$b5 = @'
THE BANK
'@
$c5 = @'
the Bank
holds the Title Deed cards and houses and hotels prior to purchase
and use by the players, pays salaries and bonuses and sells
and auctions properties and hands out their proper Title Deed cards;
it sells houses and hotels to the players and loans money when
required on mortgages.
'@
$ws.Range("B5").Value = $b5
$ws.Range("C5").Value = $c5
$ws.Range("C5").WrapText = $true
$ws.Rows(5).RowHeight = 87

# Row 6
$b6 = @'
THE PLAY
'@
$c6 = @'
Starting with the Banker, each player in turn throws the dice. The player with the highest total starts the play: Place your token on the corner marked "GO," throw the dice and move your token in the direction of the arrow the number of spaces indicated by the dice. After you have completed your play, the turn passes to the left. The tokens remain on the spaces qccupied and proceed from that point on the player's next turn. Two or more tokens may rest on the same space at the same time.
'@
$ws.Range("B6").Value = $b6
$ws.Range("C6").Value = $c6

# Row 7
$b7 = @'
"GO"
'@
$c7 = @'
Each time a player's token lands on or passes over GO, whether by throwing the dice or drawing a card, the Banker pays himther a $200 salary.
'@
$ws.Range("B7").Value = $b7
$ws.Range("C7").Value = $c7

# Row 8
$b8 = @'
BUYING PROPERTY
'@
$c8 = @'
 Whenever you land on an unowned property you
may buy that property from the Bank at its printed price. You receive the
Title Deed card showing ownership; place it faceup in front of you. 
'@
$ws.Range("B8").Value = $b8
$ws.Range("C8").Value = $c8
$ws.Range("C8").WrapText = $true
$ws.Rows(8).RowHeight = 43.5

# Row 9
$b9 = @'
PAYING RENT
'@
$c9 = @'
When you land on property owned by another player, the owner collects rent from you in accordance with the list printed on its Title Deed card.
'@
$ws.Range("B9").Value = $b9
$ws.Range("C9").Value = $c9

# Row 10
$b10 = @'
"CHANCE" AND "COMMUNITY CHEST"
'@
$c10 = @'
 When you land on either of
these spaces, take the top card from the deck indicated, follow the 
instructions and return the card facedown to the bottom of the deck. 
'@
$ws.Range("B10").Value = $b10
$ws.Range("C10").Value = $c10
$ws.Range("C10").WrapText = $true
$ws.Rows(10).RowHeight = 43.5

# Row 11
$b11 = @'
"INCOME TAX"
'@
$c11 = @'
If you land here you have two options: You may
estimate your tax at $900 and pay the Bank, or you may pay 10% of
your total worth to the Bank. Your total worth is all your cash on hand,
printed prices of mortgaged and unmortgaged properties and cost
price of all buildings you own. 
'@
$ws.Range("B11").Value = $b11
$ws.Range("C11").Value = $c11
$ws.Range("C11").WrapText = $true
$ws.Rows(11).RowHeight = 72.5

# Row 12
$b12 = @'
"JAIL"
'@
$c12 = @'
You land in Jail when. ..(I) your token lands on the space
marked "Go to Jail"; (2) you draw a card marked "Go to JailN; or
(3) you throw doubles three times in succession. 
'@
$ws.Range("B12").Value = $b12
$ws.Range("C12").Value = $c12
$ws.Range("C12").WrapText = $true
$ws.Rows(12).RowHeight = 43.5

# Row 13
$b13 = @'
"FREE PARKING
'@
$c13 = @'
A player landing on this place does not receive any money, property or reward of any kind. This is just a "free" resting place.
'@
$ws.Range("B13").Value = $b13
$ws.Range("C13").Value = $c13

# Row 14
$b14 = @'
HOUSES
'@
$c14 = @'
When you own all the propert~es in a color-group you may buy houses from the Bank and erect them on those properties.
'@
$ws.Range("B14").Value = $b14
$ws.Range("C14").Value = $c14

# Row 15
$b15 = @'
HOTELS 
'@
$c15 = @'
When a player has four houses on each property of a complete color-group, he/she may buy a hotel from the Bank and erect it on any property of the color-group. He/she returns the four houses from that property to the Bank and pays the price for the hotel as shown on the Ttle Deed card. Only one hotel may be erected on any one property
'@
$ws.Range("B15").Value = $b15
$ws.Range("C15").Value = $c15

# Row 16
$b16 = @'
BUILDING SHORTAGES
'@
$c16 = @'
When the Bank has no houses to sell, players
wishing to build must wait for some player to return or sell histher
houses to the Bank before building. If there are a limited number of
houses and hotels available and two or more players wish to buy more
than the Bank has, the houses or hotels must be sold at auction to the
highest bidder. 
'@
$ws.Range("B16").Value = $b16
$ws.Range("C16").Value = $c16
$ws.Range("C16").WrapText = $true
$ws.Rows(16).RowHeight = 87

# Row 17
$b17 = @'
SELLING PROPERTY
'@
$c17 = @'
Unimproved properties, railroads and utilities (but not buildings) may be sold to any player as a private transaction for any amount the owner can get; however, no property can be sold to - another player if buildings are standing on any properties of that colorgroup. Any buildings so located must be sold back to the Bank before the owner can sell any property of that color-group. Houses and hotels may be sold back to the Bank at any time for onehalf the price paid for them. All houses on one color-group must be sold one by one, evenly, in reverse of the manner in which they were erected. All hotels on one color-group may be sold at once, or they may be sold one house at a time (one hotel equals five houses), evenly, in reverse of the manner in which they were erected.
'@
$ws.Range("B17").Value = $b17
$ws.Range("C17").Value = $c17

# Row 18
$b18 = @'
MORTGAGES
'@
$c18 = @'
Unimproved properties can be mortgaged through the Bank at any time. Before an improved property can be mortgaged, all the buildings on all the properties of its color-group must be sold back to the Bank at half price. The mortgage value is printed on each Title Deed card. No rent can be collected on mortgaged properties or utilities, but rent can be collected on unmortgaged properties in the same group
'@
$ws.Range("B18").Value = $b18
$ws.Range("C18").Value = $c18

# Row 19
$b19 = @'
BANKRUPTCY
'@
$c19 = @'
 You are declared bankrupt if you owe more than you
can pay either to another player or to the Bank. If your ,
debt is to another player, you must tum over to that
player all that you have of value and retire from the
game
'@
$ws.Range("B19").Value = $b19
$ws.Range("C19").Value = $c19
$ws.Range("C19").WrapText = $true
$ws.Rows(19).RowHeight = 72.5

# Row 20
$b20 = @'
MISCEUANEOUS
'@
$c20 = @'
Money can be loaned to a player only by the Bank and then only by mortgaging property. No player may borrow from or lend money to another player.
'@
$ws.Range("B20").Value = $b20
$ws.Range("C20").Value = $c20

# Update the active selection to match the author's last-edited cell
$null = $ws.Range("C19").Select()

